# Updated cryptos list with latest price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.071.13"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.481.38"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.55"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.09"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.481.03"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.92"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.330"
$ws.Range("E13").Value = "  -2.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.34"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.941.92"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.460.91"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("E19").Value = "  -6.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.39"
$ws.Range("E20").Value = "  -5.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.49"
$ws.Range("E21").Value = "  -4.01%  "
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.54"
$ws.Range("E24").Value = "  -4.11%  "
$ws.Range("E25").Value = "  -5.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.78"
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.26"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("E30").Value = "  -3.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "507.02"
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("E32").Value = "  -6.06%  "
$ws.Range("E33").Value = "  -3.80%  "
$ws.Range("E34").Value = "  -4.29%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.84"
$ws.Range("E37").Value = "  -8.43%  "
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.21"
$ws.Range("E39").Value = "  -4.73%  "
$ws.Range("E40").Value = "  -6.57%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("E42").Value = "  -4.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.327"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.79"
$ws.Range("E44").Value = "  -3.76%  "
$ws.Range("E45").Value = "  -5.29%  "
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.36"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("E48").Value = "  -4.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.43"
$ws.Range("E49").Value = "  -5.33%  "
$ws.Range("E50").Value = "  -6.38%  "
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.57"
$ws.Range("E51").Value = "  -5.30%  "
